# Fixed a typo on resume
#
# 1. Update the Undergraduate Teaching Fellow date range.
# 2. Fix grammar/typos in the three existing TA bullet points.
# 3. Add a new bullet point describing exam review sessions.
# 4. Move the "_GoBack" bookmark from the end of the "Relevant Coursework"
#    block to the end of the newly added bullet point.

$d = $word.ActiveDocument

# Locate the "Undergraduate Teaching Fellow" paragraph by content instead of
# a hard-coded index, since that's more robust.
$paras = $d.Paragraphs
$tfIndex = 0
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -like "*Undergraduate Teaching Fellow*") {
        $tfIndex = $i
        break
    }
}

# 1. "Fall 2015-Present" -> "Fall 2015, 2016" (only within this paragraph,
#    so the unrelated "Association for Computing Machinery" line, which has
#    the same date text, is left untouched).
$tfRange = $paras.Item($tfIndex).Range
$tfRange.Find.Execute("Fall 2015-Present", $false, $false, $false, $false, $false, $true, 1, $false, "Fall 2015, 2016", 2) | Out-Null

# Re-resolve paragraph indices for the bullet list beneath the heading.
$bullet1Index = $tfIndex + 2
$bullet2Index = $tfIndex + 3
$bullet3Index = $tfIndex + 4

# 2a. "Increase student's understanding of course material"
#     -> "Increased students' understanding of course material"
$r1 = $paras.Item($bullet1Index).Range
$r1.Find.Execute("Increase student", $false, $false, $false, $false, $false, $true, 1, $false, "Increased students", 2) | Out-Null
$r1b = $paras.Item($bullet1Index).Range
$r1b.Find.Execute("students" + [char]0x0027 + "s", $false, $false, $false, $false, $false, $true, 1, $false, "students" + [char]0x2019, 2) | Out-Null

# 2b. "Assist the professor in leading students through class exercises"
#     -> "Assisted the professors in leading students through class exercises"
$r2 = $paras.Item($bullet2Index).Range
$r2.Find.Execute("Assist the professor ", $false, $false, $false, $false, $false, $true, 1, $false, "Assisted the professors ", 2) | Out-Null

# 2c. "Co-lead weekly study hall and tutoring sessions..."
#     -> "Co-lead weekly study halls and tutoring sessions..."
$r3 = $paras.Item($bullet3Index).Range
$r3.Find.Execute("study hall and", $false, $false, $false, $false, $false, $true, 1, $false, "study halls and", 2) | Out-Null

# 3. Add a new bullet point after the "Co-lead weekly study halls..." line.
$bulletPara = $paras.Item($bullet3Index)
$bulletPara.Range.InsertParagraphAfter()

$newBulletIndex = $bullet3Index + 1
$newPara = $d.Paragraphs.Item($newBulletIndex)
$newRange = $newPara.Range
$newRange.Text = "Created and hosted exam review sessions to facilitate student success"

# 4. Relocate the "_GoBack" bookmark to the end of the text we just typed
#    (matching real Word's behaviour of tracking the most recent edit).
try {
    $oldGoBack = $d.Bookmarks.Item("_GoBack")
    $oldGoBack.Delete()
} catch {
}

$finalPara = $d.Paragraphs.Item($newBulletIndex)
$finalRange = $finalPara.Range
$finalRange.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $finalRange) | Out-Null
